$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ40571502",
    "summ40653499",
    "summ40755704",
    "summ40854826",
    "summ40951766",
    "summ41037521",
    "summ41137559",
    "summ41237593",
    "summ41321642",
    "summ41420581",
    "summ41527673",
    "summ41626098",
    "summ41732174",
    "summ41830198",
    "summ41932709",
    "summ42049018",
    "summ42159671",
    "summ42255494",
    "summ42354796",
    "summ42461298",
    "summ42557815",
    "summ42657455",
    "summ42763444",
    "summ42865962",
    "summ42964616",
    "summ43064373",
    "summ43170004",
    "summ43292318",
    "summ43397653",
    "summ43497165",
    "summ43594972",
    "summ43698005",
    "summ43796025",
    "summ43895540",
    "summ44111677",
    "summ44211313",
    "summ44313622",
    "summ44454080",
    "summ44555119",
    "summ44654278",
    "summ44738339",
    "summ44840613",
    "summ44936773",
    "summ45037830",
    "summ45121648",
    "summ45220900",
    "summ45320853",
    "summ45421411",
    "summ45520302",
    "summ45605896"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

